$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow writing the refreshed
# holdings data, then restore protection afterwards.
$ws.Unprotect("D382")

$ws.Range("D2").Value = 0.08405998806056415
$ws.Range("E2").Value = 0.0179552175749893
$ws.Range("D3").Value = 0.05020679544274447
$ws.Range("E3").Value = 0.005669606664236326
$ws.Range("D4").Value = 0.04298220584784766
$ws.Range("E4").Value = 0.008175466397752684
$ws.Range("D5").Value = 0.03599956201170304
$ws.Range("E5").Value = 0.001669449081802998
$ws.Range("D6").Value = 0.03447248998525426
$ws.Range("E6").Value = 0.002250574530469951
$ws.Range("D7").Value = 0.02982584747472108
$ws.Range("E7").Value = 0.008572959457045837
$ws.Range("D8").Value = 0.02724047902064036
$ws.Range("E8").Value = 0.00117332638028822
$ws.Range("D9").Value = 0.02777795920168472
$ws.Range("E9").Value = 0.006923409779316492
$ws.Range("D10").Value = 0.02662195566109037
$ws.Range("E10").Value = 0.007501103103397577
$ws.Range("D11").Value = 0.02691070613192335
$ws.Range("E11").Value = 0.004601932811780785
$ws.Range("D12").Value = 0.02225746178943324
$ws.Range("E12").Value = -0.00798602445719987
$ws.Range("D13").Value = 0.02214263544325841
$ws.Range("E13").Value = -0.005008923944959331
$ws.Range("D14").Value = 0.02026225158242554
$ws.Range("E14").Value = 0.005285013212533052
$ws.Range("D15").Value = 0.02054699542420881
$ws.Range("E15").Value = 0.003412470140886414
$ws.Range("D16").Value = 0.02183248593090991
$ws.Range("E16").Value = 0.0002252252252250564
$ws.Range("D17").Value = 0.01892003085077051
$ws.Range("E17").Value = -0.01386108116433094
$ws.Range("D18").Value = 0.01784006220236964
$ws.Range("E18").Value = -0.001643561526368353
$ws.Range("D19").Value = 0.01712469679749298
$ws.Range("E19").Value = -0.005849197064766498
$ws.Range("D20").Value = 0.01766823245198748
$ws.Range("E20").Value = -0.006491297692612807
$ws.Range("D21").Value = 0.01508850970247681
$ws.Range("E21").Value = -0.01564282222584334
$ws.Range("D22").Value = 0.01453586807286932
$ws.Range("E22").Value = 0.0132180667794275
$ws.Range("D23").Value = 0.01347702528673059
$ws.Range("E23").Value = -0.003101309441764233
$ws.Range("D24").Value = 0.01334699196211707
$ws.Range("E24").Value = -0.00878054770968939
$ws.Range("D25").Value = 0.01197162548105074
$ws.Range("E25").Value = -0.001977637483836503
$ws.Range("D26").Value = 0.01104190542201424
$ws.Range("E26").Value = 0.005690252350321501
$ws.Range("D27").Value = 0.01181363681284032
$ws.Range("E27").Value = 0.008844953173777315
$ws.Range("D28").Value = 0.01030322872084832
$ws.Range("E28").Value = 0.008166295471417895
$ws.Range("D29").Value = 0.01033637447025961
$ws.Range("E29").Value = -0.03527380365071531
$ws.Range("D30").Value = 0.01013904798956107
$ws.Range("E30").Value = 0.01131617944227403
$ws.Range("D31").Value = 0.01027336112247771
$ws.Range("E31").Value = 0.002127282396738206
$ws.Range("D32").Value = 0.01021508288175456
$ws.Range("E32").Value = -0.002317703690497552
$ws.Range("D33").Value = 0.0102130795672297
$ws.Range("E33").Value = 0.03623013962445865
$ws.Range("D34").Value = 0.0102237335581119
$ws.Range("E34").Value = 0.01656646626586511
$ws.Range("D35").Value = 0.01001766534130489
$ws.Range("E35").Value = 0.00152710613387641
$ws.Range("D36").Value = 0.01093468256505876
$ws.Range("E36").Value = 0.0003122853038535212
$ws.Range("D37").Value = 0.008156039789204782
$ws.Range("E37").Value = -0.006698821007502631
$ws.Range("D38").Value = 0.009602114167023492
$ws.Range("E38").Value = 0.005761106132377458
$ws.Range("D39").Value = 0.008494508884154697
$ws.Range("E39").Value = -0.01336227689339109
$ws.Range("D40").Value = 0.009737747666331511
$ws.Range("E40").Value = 0.008836896145430151
$ws.Range("D41").Value = 0.008522373168000454
$ws.Range("E41").Value = -0.01019328781613604
$ws.Range("D42").Value = 0.008826831445903346
$ws.Range("E42").Value = 0.004425668879501021
$ws.Range("D43").Value = 0.008957593248525911
$ws.Range("E43").Value = 0.008813617903741466
$ws.Range("D44").Value = 0.008131453656399702
$ws.Range("E44").Value = 0.001970928800197003
$ws.Range("D45").Value = 0.0089997083834235
$ws.Range("E45").Value = -0.008640838586302158
$ws.Range("D46").Value = 0.008025642225586733
$ws.Range("E46").Value = 0.0007828810020877874
$ws.Range("D47").Value = 0.009416670983847408
$ws.Range("E47").Value = 0.01094650524117147
$ws.Range("D48").Value = 0.008128721863865804
$ws.Range("E48").Value = 0.0092530358023033
$ws.Range("D49").Value = 0.009009178597541013
$ws.Range("E49").Value = -0.02263561660450608
$ws.Range("D50").Value = 0.00680507732144127
$ws.Range("E50").Value = 0.005479580367178549
$ws.Range("D51").Value = 0.007696005926496418
$ws.Range("E51").Value = 0.01405651000993902
$ws.Range("D52").Value = 0.007903758748699334
$ws.Range("E52").Value = 0.01218353063164246
$ws.Range("D53").Value = 0.008574914644402418
$ws.Range("E53").Value = 0.001699090986322371
$ws.Range("D54").Value = 0.006609414957694629
$ws.Range("E54").Value = 0.003846960611793904
$ws.Range("D55").Value = 0.006464513852215364
$ws.Range("E55").Value = -0.007902298850574585
$ws.Range("D56").Value = 0.005590067062114732
$ws.Range("E56").Value = -0.004104969945755488
$ws.Range("D57").Value = 0.006156367654391712
$ws.Range("E57").Value = -0.0003549875754348397
$ws.Range("D58").Value = 0.005906135458286688
$ws.Range("E58").Value = -0.01618871415356149
$ws.Range("D59").Value = 0.006207998533282377
$ws.Range("E59").Value = 0.02299229922992296
$ws.Range("D60").Value = 0.005148063030130093
$ws.Range("E60").Value = 0.009286282833643078
$ws.Range("D61").Value = 0.005012793769826596
$ws.Range("E61").Value = 0.003097212508742286
$ws.Range("D62").Value = 0.005313518597933161
$ws.Range("E62").Value = 0.00548395941870039
$ws.Range("D63").Value = 0.004855305930247398
$ws.Range("E63").Value = 0.007351837959489904
$ws.Range("D64").Value = 0.004904113956853035
$ws.Range("E64").Value = 0.008021390374331583
$ws.Range("D65").Value = 0.00454114978884917
$ws.Range("E65").Value = 0.01054742330058156
$ws.Range("D66").Value = 0.004282631155391323
$ws.Range("E66").Value = 0.007973464310773837
$ws.Range("D67").Value = 0.003990056175010886
$ws.Range("E67").Value = 0.004929481035190975
$ws.Range("D68").Value = 0.004024066992057912
$ws.Range("E68").Value = -0.003597976986524531
$ws.Range("D69").Value = 0.004130379251502094
$ws.Range("E69").Value = 0.008576026808351234
$ws.Range("D70").Value = 0.00399715883559902
$ws.Range("E70").Value = 0.015343083652269
$ws.Range("D71").Value = 0.003351545200087879
$ws.Range("E71").Value = -0.00896592946802155
$ws.Range("D72").Value = 0.003418155408039416
$ws.Range("E72").Value = -0.007232767232767334
$ws.Range("D73").Value = 0.003769600517525347
$ws.Range("E73").Value = -0.01578616807981259
$ws.Range("D74").Value = 0.002717678272472498
$ws.Range("E74").Value = 0.01486011057128489
$ws.Range("D75").Value = 0.002466080180100525
$ws.Range("E75").Value = -0.01048667011299009
$ws.Range("D76").Value = 0.001935383950515345
$ws.Range("E76").Value = -0.0219252846523007
$ws.Range("D77").Value = 0.00159482048128944
$ws.Range("E77").Value = -0.0375699440447641
$ws.Range("E78").Value = 0.003130244963410522

$ws.Protect("D382")
